$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 317.85715
$ws.Range("I4").Value = 321.54544
$ws.Range("J4").Value = 304.33334
$ws.Range("K4").Value = 321.54544
$ws.Range("L4").Value = 304.33334
$ws.Range("M4").Value = -207.54544
$ws.Range("N4").Value = -532.33334

$ws.Range("H38").Value = 430.8
$ws.Range("I38").Value = 51
$ws.Range("K38").Value = 153
$ws.Range("M38").Value = 219

$ws.Range("H43").Value = 11218.5
$ws.Range("I43").Value = 7999
$ws.Range("K43").Value = 7999
$ws.Range("M43").Value = -7930

$ws.Range("H88").Value = 2388.5
$ws.Range("I88").Value = 3081.3333
$ws.Range("J88").Value = 2249.9333
$ws.Range("K88").Value = 3081.3333
$ws.Range("L88").Value = 2249.9333
$ws.Range("M88").Value = -2675.3333
$ws.Range("N88").Value = -3061.9333

$ws.Range("H91").Value = 2388.5
$ws.Range("I91").Value = 3081.3333
$ws.Range("J91").Value = 2249.9333
$ws.Range("K91").Value = 3081.3333
$ws.Range("L91").Value = 2249.9333
$ws.Range("M91").Value = -1677.3333
$ws.Range("N91").Value = -5057.933300000001

$ws.Range("H113").Value = 27246
$ws.Range("I113").Value = 29161.555
$ws.Range("K113").Value = 29161.555
$ws.Range("M113").Value = -25907.555

$ws.Range("H125").Value = 6268
$ws.Range("I125").Value = 3500
$ws.Range("J125").Value = 9036
$ws.Range("K125").Value = 31500
$ws.Range("L125").Value = 81324
$ws.Range("M125").Value = -29040
$ws.Range("N125").Value = -86244

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3183.1667
$ws.Range("I32").Value = 2335.25
$ws.Range("J32").Value = 7422.75
$ws.Range("K32").Value = 2335.25
$ws.Range("L32").Value = 7422.75
$ws.Range("M32").Value = -2048.25
$ws.Range("N32").Value = -7996.75

$ws.Range("H132").Value = 311.66666
$ws.Range("I132").Value = 311.66666
$ws.Range("K132").Value = 934.9999799999999
$ws.Range("M132").Value = 1595.00002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 19000000
$ws.Range("I7").Value = 19000000
$ws.Range("K7").Value = 19000000
$ws.Range("M7").Value = -18999887

$ws.Range("H22").Value = 3499.75
$ws.Range("I22").Value = 2999.6667
$ws.Range("J22").Value = 5000
$ws.Range("K22").Value = 2999.6667
$ws.Range("L22").Value = 5000
$ws.Range("M22").Value = -2826.6667
$ws.Range("N22").Value = -5346

$ws.Range("H64").Value = 910.3333
$ws.Range("I64").Value = 886.75
$ws.Range("K64").Value = 886.75
$ws.Range("M64").Value = -661.75

$ws.Range("H67").Value = 910.3333
$ws.Range("I67").Value = 886.75
$ws.Range("K67").Value = 886.75
$ws.Range("M67").Value = -106.75

$ws.Range("H80").Value = 343.36365
$ws.Range("I80").Value = 546.75
$ws.Range("J80").Value = 227.14285
$ws.Range("K80").Value = 546.75
$ws.Range("L80").Value = 227.14285
$ws.Range("M80").Value = 451.25
$ws.Range("N80").Value = -2223.14285

$ws.Range("H83").Value = 343.36365
$ws.Range("I83").Value = 546.75
$ws.Range("J83").Value = 227.14285
$ws.Range("K83").Value = 2733.75
$ws.Range("L83").Value = 1135.71425
$ws.Range("M83").Value = 2258.25
$ws.Range("N83").Value = -11119.71425

$ws.Range("H86").Value = 4017.5454
$ws.Range("I86").Value = 813.2857
$ws.Range("K86").Value = 813.2857
$ws.Range("M86").Value = 309.7143

$ws.Range("H89").Value = 4017.5454
$ws.Range("I89").Value = 813.2857
$ws.Range("K89").Value = 4066.4285
$ws.Range("M89").Value = 1549.5715

$ws.Range("H99").Value = 4967
$ws.Range("I99").Value = 4646
$ws.Range("J99").Value = 5159.6
$ws.Range("K99").Value = 4646
$ws.Range("L99").Value = 5159.6
$ws.Range("M99").Value = -3148
$ws.Range("N99").Value = -8155.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 25313.5
$ws.Range("I12").Value = 627
$ws.Range("J12").Value = 50000
$ws.Range("K12").Value = 627
$ws.Range("L12").Value = 50000
$ws.Range("M12").Value = -457
$ws.Range("N12").Value = -50340

$ws.Range("H58").Value = 1509.4375
$ws.Range("I58").Value = 1096
$ws.Range("K58").Value = 1096
$ws.Range("M58").Value = -893

$ws.Range("H105").Value = 6144.15
$ws.Range("I105").Value = 7333.375
$ws.Range("K105").Value = 7333.375
$ws.Range("M105").Value = -5586.375

$ws.Range("H132").Value = 2457.0557
$ws.Range("I132").Value = 2366.9355
$ws.Range("K132").Value = 7100.806500000001
$ws.Range("M132").Value = -4570.806500000001

$ws.Range("H136").Value = 1509.4375
$ws.Range("I136").Value = 1096
$ws.Range("K136").Value = 3288
$ws.Range("M136").Value = -738

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 266.1111
$ws.Range("J23").Value = 321.5
$ws.Range("L23").Value = 964.5
$ws.Range("N23").Value = -1434.5

$ws.Range("H68").Value = 3000
$ws.Range("J68").Value = 3000
$ws.Range("L68").Value = 9000
$ws.Range("N68").Value = -10622

$ws.Range("H71").Value = 3000
$ws.Range("J71").Value = 3000
$ws.Range("L71").Value = 27000
$ws.Range("N71").Value = -35112

$ws.Range("H122").Value = 509.57144
$ws.Range("I122").Value = 135
$ws.Range("J122").Value = 1009
$ws.Range("K122").Value = 1215
$ws.Range("L122").Value = 9081
$ws.Range("M122").Value = 1235
$ws.Range("N122").Value = -13981

$ws.Range("H131").Value = 1134.5714
$ws.Range("I131").Value = 993.5
$ws.Range("J131").Value = 1158.0834
$ws.Range("K131").Value = 2980.5
$ws.Range("L131").Value = 3474.2502
$ws.Range("M131").Value = 2059.5
$ws.Range("N131").Value = -13554.2502

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1374.1428
$ws.Range("I102").Value = 840.0909
$ws.Range("J102").Value = 3332.3333
$ws.Range("K102").Value = 840.0909
$ws.Range("L102").Value = 3332.3333
$ws.Range("M102").Value = 781.9091
$ws.Range("N102").Value = -6576.3333

$ws.Range("H107").Value = 14072
$ws.Range("I107").Value = 2296
$ws.Range("J107").Value = 25848
$ws.Range("K107").Value = 2296
$ws.Range("L107").Value = 25848
$ws.Range("M107").Value = -376
$ws.Range("N107").Value = -29688

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 471.5
$ws.Range("I16").Value = 471.5
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 471.5
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -301.5
$ws.Range("N16").ClearContents()

$ws.Range("H55").Value = 418.46155
$ws.Range("I55").Value = 132.75
$ws.Range("J55").Value = 545.44446
$ws.Range("K55").Value = 132.75
$ws.Range("L55").Value = 545.44446
$ws.Range("M55").Value = 40.25
$ws.Range("N55").Value = -891.44446

$ws.Range("H93").Value = 1999.3334
$ws.Range("I93").Value = 1749.75
$ws.Range("J93").Value = 2498.5
$ws.Range("K93").Value = 1749.75
$ws.Range("L93").Value = 2498.5
$ws.Range("M93").Value = -501.75
$ws.Range("N93").Value = -4994.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2343.75
$ws.Range("I81").Value = 2625
$ws.Range("J81").Value = 1500
$ws.Range("K81").Value = 5250
$ws.Range("L81").Value = 3000
$ws.Range("M81").Value = -4189
$ws.Range("N81").Value = -5122

$ws.Range("H84").Value = 2343.75
$ws.Range("I84").Value = 2625
$ws.Range("J84").Value = 1500
$ws.Range("K84").Value = 26250
$ws.Range("L84").Value = 15000
$ws.Range("M84").Value = -20946
$ws.Range("N84").Value = -25608
